$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.538.75'
Set-TextValue $ws.Range('E2') '  +3.37%  '

Set-TextValue $ws.Range('D3') '1.824.25'
Set-TextValue $ws.Range('E3') '  +4.51%  '

Set-TextValue $ws.Range('D4') '1.004'
Set-TextValue $ws.Range('E4') '  +0.04%  '

Set-TextValue $ws.Range('D5') '343.75'
Set-TextValue $ws.Range('E5') '  +2.61%  '

Set-TextValue $ws.Range('D6') '1.002'
Set-TextValue $ws.Range('E6') '  +0.07%  '

Set-TextValue $ws.Range('D7') '0.3822'
Set-TextValue $ws.Range('E7') '  +0.38%  '

Set-TextValue $ws.Range('D8') '0.3534'
Set-TextValue $ws.Range('E8') '  +4.20%  '

Set-TextValue $ws.Range('D9') '49.87'
Set-TextValue $ws.Range('E9') '  +1.92%  '

Set-TextValue $ws.Range('D10') '1.238'
Set-TextValue $ws.Range('E10') '  +3.57%  '

Set-TextValue $ws.Range('D11') '0.07739'
Set-TextValue $ws.Range('E11') '  +3.40%  '

Set-TextValue $ws.Range('D12') '1.002'
Set-TextValue $ws.Range('E12') '  +0.00%  '

Set-TextValue $ws.Range('D13') '22.15'
Set-TextValue $ws.Range('E13') '  +7.97%  '

Set-TextValue $ws.Range('D14') '6.609'
Set-TextValue $ws.Range('E14') '  +2.05%  '

Set-TextValue $ws.Range('D15') '1.830.21'
Set-TextValue $ws.Range('E15') '  +4.86%  '

Set-TextValue $ws.Range('D16') '7.214'
Set-TextValue $ws.Range('E16') '  +1.25%  '

Set-TextValue $ws.Range('D17') '0.00001126'
Set-TextValue $ws.Range('E17') '  +3.26%  '

Set-TextValue $ws.Range('D18') '0.06745'
Set-TextValue $ws.Range('E18') '  +0.55%  '

Set-TextValue $ws.Range('D19') '86.93'
Set-TextValue $ws.Range('E19') '  +3.97%  '

Set-TextValue $ws.Range('E20') '  +0.09%  '

Set-TextValue $ws.Range('D21') '17.59'
Set-TextValue $ws.Range('E21') '  +4.66%  '

Set-TextValue $ws.Range('D22') '6.544'
Set-TextValue $ws.Range('E22') '  +5.23%  '

Set-TextValue $ws.Range('D23') '13.19'
Set-TextValue $ws.Range('E23') '  +0.69%  '

Set-TextValue $ws.Range('D24') '27.515.85'
Set-TextValue $ws.Range('E24') '  +3.41%  '

Set-TextValue $ws.Range('D25') '2.478'
Set-TextValue $ws.Range('E25') '  +1.20%  '

Set-TextValue $ws.Range('D26') '2.687'
Set-TextValue $ws.Range('E26') '  +8.44%  '

Set-TextValue $ws.Range('D27') '22.03'
Set-TextValue $ws.Range('E27') '  +11.81%  '

Set-TextValue $ws.Range('D28') '1.484'
Set-TextValue $ws.Range('E28') '  +3.45%  '

Set-TextValue $ws.Range('D29') '153.08'
Set-TextValue $ws.Range('E29') '  -0.34%  '

Set-TextValue $ws.Range('D30') '2.034.99'
Set-TextValue $ws.Range('E30') '  +5.04%  '

Set-TextValue $ws.Range('D31') '135.36'
Set-TextValue $ws.Range('E31') '  +2.39%  '

Set-TextValue $ws.Range('D32') '6.345'
Set-TextValue $ws.Range('E32') '  +3.33%  '

Set-TextValue $ws.Range('D33') '4.096'
Set-TextValue $ws.Range('E33') '  -0.97%  '

Set-TextValue $ws.Range('D34') '13.89'
Set-TextValue $ws.Range('E34') '  +6.44%  '

Set-TextValue $ws.Range('D35') '0.08797'
Set-TextValue $ws.Range('E35') '  +1.17%  '

Set-TextValue $ws.Range('E36') '  -0.81%  '

Set-TextValue $ws.Range('D37') '5.625'
Set-TextValue $ws.Range('E37') '  +3.01%  '

Set-TextValue $ws.Range('D38') '0.7005'
Set-TextValue $ws.Range('E38') '  +11.77%  '

Set-TextValue $ws.Range('D39') '9.112'
Set-TextValue $ws.Range('E39') '  +5.46%  '

Set-TextValue $ws.Range('D40') '0.06530'
Set-TextValue $ws.Range('E40') '  +2.90%  '

Set-TextValue $ws.Range('D41') '0.2259'
Set-TextValue $ws.Range('E41') '  +3.31%  '

Set-TextValue $ws.Range('D42') '0.02405'
Set-TextValue $ws.Range('E42') '  +1.49%  '

Set-TextValue $ws.Range('D43') '1.307'
Set-TextValue $ws.Range('E43') '  +6.20%  '

Set-TextValue $ws.Range('D44') '14.68'
Set-TextValue $ws.Range('E44') '  +2.63%  '

Set-TextValue $ws.Range('D45') '0.6600'
Set-TextValue $ws.Range('E45') '  +8.47%  '

Set-TextValue $ws.Range('D46') '1.002'
Set-TextValue $ws.Range('E46') '  +0.04%  '

Set-TextValue $ws.Range('D47') '3.950'
Set-TextValue $ws.Range('E47') '  +0.42%  '

Set-TextValue $ws.Range('D48') '2.191'
Set-TextValue $ws.Range('E48') '  +5.61%  '

Set-TextValue $ws.Range('D49') '133.19'
Set-TextValue $ws.Range('E49') '  +3.25%  '

Set-TextValue $ws.Range('D50') '0.07306'
Set-TextValue $ws.Range('E50') '  +0.57%  '

Set-TextValue $ws.Range('D51') '81.11'
Set-TextValue $ws.Range('E51') '  +3.77%  '
